$d = $word.ActiveDocument

function Find-ParagraphIndex($startsWith) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($startsWith)) {
            return $i
        }
    }
    return -1
}

# --- Edit 1: paragraph "6. ... two million elements." loses its yellow highlighting ---
$idx6 = Find-ParagraphIndex("6.")
$p6 = $d.Paragraphs.Item($idx6)
$xml6 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="707D6A73" w14:textId="77777777" w:rsidR="004348C3" w:rsidRPr="00D91C12" w:rsidRDefault="00136598"><w:pPr><w:ind w:left="720" w:hanging="720"/></w:pPr><w:r w:rsidRPr="00D91C12"><w:lastRenderedPageBreak/><w:t xml:space="preserve">6. </w:t></w:r><w:r w:rsidRPr="00D91C12"><w:tab/><w:t>Now analyze the runtimes you recorded and give an estimate of how long you believe it would take each algorithm to sort a vector containing two million elements.</w:t></w:r></w:p>'
$p6.Range.InsertXML($xml6)

# --- Edit 2: "Turn In:" block is reorganized -
#     the trailing bookmark paragraph merges into "Turn In:" (keeping the bookmark
#     ahead of the text), the "Submit..." and "spreadsheet..." paragraphs lose their
#     yellow highlighting, and a new empty highlighted paragraph is appended ---
$idxTurnIn = Find-ParagraphIndex("Turn In")
$idxNote = Find-ParagraphIndex("Note:")
$idxBookmarkPara = $idxNote - 1

$pStart = $d.Paragraphs.Item($idxTurnIn)
$pEnd = $d.Paragraphs.Item($idxBookmarkPara)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="448FE82C" w14:textId="77777777" w:rsidR="004348C3" w:rsidRPr="00D91C12" w:rsidRDefault="00136598"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Turn </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>In</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Submit your code and spreadsheet through the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Classroom repository. I should be able to see you generating each vector for each algorithm (you may comment these out so you’re not testing every algorithm every time you run your program).</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr></w:p><w:p><w:r><w:t xml:space="preserve">The spreadsheet should contain your record of </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>sizes</w:t></w:r><w:r><w:t xml:space="preserve"> needed for a 4.0 second base sort time, and the </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>times</w:t></w:r><w:r><w:t xml:space="preserve"> for doubled and quadrupled sizes. It should also contain your </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>estimates</w:t></w:r><w:r><w:t xml:space="preserve"> of the time required to sort two million numbers for the 5 selected combinations.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr></w:p>'
$rng.InsertXML($xml2)

Write-Host "Edits applied successfully"
